# Horarios Línea 141 - actualización 04:40:32
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173):
#  - refresh the "Última actualización" timestamp and "Total filas" count
#  - shift the tail of each sheet's data down to make room for newly
#    scraped rows, then fill in the new rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:40:32"
$ws1.Range("A3").Value = "Total filas: 28"

$data1 = @(
    @("04:40:32", "05:17", "14_ABASTO", 37, "LP1912"),
    @("03:23:38", "05:22", "23_HERNANDEZ", 119, "LP1912"),
    @("03:54:40", "05:34", "215B_EL PATO", 100, "LP1912"),
    @("04:18:06", "05:35", "215B_EL PATO", 77, "LP1912"),
    @("04:18:06", "05:36", "14_ABASTO", 78, "LP1912"),
    @("03:54:40", "05:39", "14_ABASTO", 105, "LP1912"),
    @("03:54:40", "05:46", "15_ABASTO", 112, "LP1912"),
    @("04:40:32", "06:04", "16_SANTA ANA", 84, "LP1912"),
    @("04:18:06", "06:09", "16_SANTA ANA", 111, "LP1912"),
    @("04:40:32", "06:11", "215A_EL PATO", 91, "LP1912"),
    @("04:18:06", "06:12", "215A_EL PATO", 114, "LP1912"),
    @("04:18:06", "06:14", "225_HARAS DEL SUR", 116, "LP1912"),
    @("04:40:32", "06:21", "26_HERNANDEZ", 101, "LP1912"),
    @("04:40:32", "06:27", "23_HERNANDEZ", 107, "LP1912"),
    @("04:40:32", "06:29", "86_EST CHICA-ESC AGRARIA", 109, "LP1912"),
    @("04:40:32", "06:31", "16_SANTA ANA", 111, "LP1912")
)

$startRow1 = 18
for ($i = 0; $i -lt $data1.Count; $i++) {
    $r = $startRow1 + $i
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:40:32"
$ws2.Range("A3").Value = "Total filas: 9"

# Row 13 becomes the newly scraped entry; the former row 13 moves to row 14
$ws2.Cells.Item(14, 1).Value = "04:18:06"
$ws2.Cells.Item(14, 2).Value = "06:12"
$ws2.Cells.Item(14, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(14, 4).Value = 114
$ws2.Cells.Item(14, 5).Value = "LP1912"

$ws2.Cells.Item(13, 1).Value = "04:40:32"
$ws2.Cells.Item(13, 2).Value = "06:11"
$ws2.Cells.Item(13, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(13, 4).Value = 91
$ws2.Cells.Item(13, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:40:32"
$ws3.Range("A3").Value = "Total filas: 5"

$ws3.Cells.Item(10, 1).Value = "04:40:32"
$ws3.Cells.Item(10, 2).Value = "06:33"
$ws3.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 113
$ws3.Cells.Item(10, 5).Value = "L6203"
